# Update stock data: adjust product names and creation_date timestamps
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductData")

# Row 9: Jabon piel dove en barra x 5 -> x 10
$ws.Range("C9").Value = "Jabon piel dove en barra x 10"
$ws.Range("G9").Value = 45813.87396958334

# Row 10: Dolex Gripa x 12 pastillas -> x 20 pastillas
$ws.Range("C10").Value = "Dolex Gripa x 20 pastillas"
$ws.Range("G10").Value = 45813.87414978183
